$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 11
$ws.Range("B2").Value = 'Vega Monumental Concepción'
$ws.Range("C2").Value = 'Bíobío'
$ws.Range("D2").Value = 45106
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 'Fruta'
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = 'Otros'
$ws.Range("I2").Value = 100107001
$ws.Range("J2").Value = 'Caqui'
$ws.Range("K2").Value = 'Mankaki'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 18000
$ws.Range("P2").Value = 17625
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = 'Región del Maule'
$ws.Range("S2").Value = 979
$ws.Range("T2").Value = 18

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = 'Vega Monumental Concepción'
$ws.Range("C3").Value = 'Bíobío'
$ws.Range("D3").Value = 45084
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 'Fruta'
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = 'Otros'
$ws.Range("I3").Value = 100107001
$ws.Range("J3").Value = 'Caqui'
$ws.Range("K3").Value = 'Mankaki'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 18000
$ws.Range("P3").Value = 17500
$ws.Range("Q3").Value = '$/caja 18 kilos granel'
$ws.Range("R3").Value = 'Región del Maule'
$ws.Range("S3").Value = 972
$ws.Range("T3").Value = 18

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 'Vega Monumental Concepción'
$ws.Range("C4").Value = 'Bíobío'
$ws.Range("D4").Value = 44742
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 'Fruta'
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = 'Otros'
$ws.Range("I4").Value = 100107001
$ws.Range("J4").Value = 'Caqui'
$ws.Range("K4").Value = 'Mankaki'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 806
$ws.Range("T4").Value = 18

$ws.Range("A5").Value = 11
$ws.Range("B5").Value = 'Vega Monumental Concepción'
$ws.Range("C5").Value = 'Bíobío'
$ws.Range("D5").Value = 45093
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 'Fruta'
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = 'Otros'
$ws.Range("I5").Value = 100107001
$ws.Range("J5").Value = 'Caqui'
$ws.Range("K5").Value = 'Mankaki'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17429
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = 'Provincia de Curicó'
$ws.Range("S5").Value = 968
$ws.Range("T5").Value = 18

$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 'Vega Monumental Concepción'
$ws.Range("C6").Value = 'Bíobío'
$ws.Range("D6").Value = 44707
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 'Fruta'
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = 'Otros'
$ws.Range("I6").Value = 100107001
$ws.Range("J6").Value = 'Caqui'
$ws.Range("K6").Value = 'Mankaki'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12500
$ws.Range("Q6").Value = '$/caja 12 kilos empedrada'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 1042
$ws.Range("T6").Value = 12

$ws.Range("A7").Value = 11
$ws.Range("B7").Value = 'Vega Monumental Concepción'
$ws.Range("C7").Value = 'Bíobío'
$ws.Range("D7").Value = 45092
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 'Fruta'
$ws.Range("G7").Value = 100107
$ws.Range("H7").Value = 'Otros'
$ws.Range("I7").Value = 100107001
$ws.Range("J7").Value = 'Caqui'
$ws.Range("K7").Value = 'Mankaki'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 140
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 19000
$ws.Range("P7").Value = 18429
$ws.Range("Q7").Value = '$/caja 18 kilos granel'
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 1024
$ws.Range("T7").Value = 18

$ws.Range("A8").Value = 11
$ws.Range("B8").Value = 'Vega Monumental Concepción'
$ws.Range("C8").Value = 'Bíobío'
$ws.Range("D8").Value = 44714
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 'Fruta'
$ws.Range("G8").Value = 100107
$ws.Range("H8").Value = 'Otros'
$ws.Range("I8").Value = 100107001
$ws.Range("J8").Value = 'Caqui'
$ws.Range("K8").Value = 'Mankaki'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("Q8").Value = '$/caja 18 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 806
$ws.Range("T8").Value = 18

$ws.Range("A9").Value = 11
$ws.Range("B9").Value = 'Vega Monumental Concepción'
$ws.Range("C9").Value = 'Bíobío'
$ws.Range("D9").Value = 44719
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 'Fruta'
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = 'Otros'
$ws.Range("I9").Value = 100107001
$ws.Range("J9").Value = 'Caqui'
$ws.Range("K9").Value = 'Mankaki'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14400
$ws.Range("Q9").Value = '$/caja 18 kilos granel'
$ws.Range("R9").Value = 'Región del Maule'
$ws.Range("S9").Value = 800
$ws.Range("T9").Value = 18

$ws.Range("A10").Value = 11
$ws.Range("B10").Value = 'Vega Monumental Concepción'
$ws.Range("C10").Value = 'Bíobío'
$ws.Range("D10").Value = 45100
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 'Fruta'
$ws.Range("G10").Value = 100107
$ws.Range("H10").Value = 'Otros'
$ws.Range("I10").Value = 100107001
$ws.Range("J10").Value = 'Caqui'
$ws.Range("K10").Value = 'Mankaki'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1000
$ws.Range("T10").Value = 18

$ws.Range("A11").Value = 11
$ws.Range("B11").Value = 'Vega Monumental Concepción'
$ws.Range("C11").Value = 'Bíobío'
$ws.Range("D11").Value = 44334
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 'Fruta'
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = 'Otros'
$ws.Range("I11").Value = 100107001
$ws.Range("J11").Value = 'Caqui'
$ws.Range("K11").Value = 'Mankaki'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 100
$ws.Range("N11").Value = 11000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 11500
$ws.Range("Q11").Value = '$/caja 12 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 11500
$ws.Range("T11").Value = 1

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 'Vega Monumental Concepción'
$ws.Range("C12").Value = 'Bíobío'
$ws.Range("D12").Value = 44330
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 'Fruta'
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = 'Otros'
$ws.Range("I12").Value = 100107001
$ws.Range("J12").Value = 'Caqui'
$ws.Range("K12").Value = 'Mankaki'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 15500
$ws.Range("Q12").Value = '$/caja 18 kilos granel'
$ws.Range("R12").Value = 'Provincia de Curicó'
$ws.Range("S12").Value = 861
$ws.Range("T12").Value = 18

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 'Vega Monumental Concepción'
$ws.Range("C13").Value = 'Bíobío'
$ws.Range("D13").Value = 44708
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 'Fruta'
$ws.Range("G13").Value = 100107
$ws.Range("H13").Value = 'Otros'
$ws.Range("I13").Value = 100107001
$ws.Range("J13").Value = 'Caqui'
$ws.Range("K13").Value = 'Mankaki'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 13000
$ws.Range("P13").Value = 12571
$ws.Range("Q13").Value = '$/caja 12 kilos empedrada'
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("S13").Value = 1048
$ws.Range("T13").Value = 12

$ws.Range("A14").Value = 11
$ws.Range("B14").Value = 'Vega Monumental Concepción'
$ws.Range("C14").Value = 'Bíobío'
$ws.Range("D14").Value = 45090
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 'Fruta'
$ws.Range("G14").Value = 100107
$ws.Range("H14").Value = 'Otros'
$ws.Range("I14").Value = 100107001
$ws.Range("J14").Value = 'Caqui'
$ws.Range("K14").Value = 'Mankaki'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 150
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17533
$ws.Range("Q14").Value = '$/caja 18 kilos granel'
$ws.Range("R14").Value = 'Región del Maule'
$ws.Range("S14").Value = 974
$ws.Range("T14").Value = 18

$ws.Range("A15").Value = 11
$ws.Range("B15").Value = 'Vega Monumental Concepción'
$ws.Range("C15").Value = 'Bíobío'
$ws.Range("D15").Value = 45090
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 'Fruta'
$ws.Range("G15").Value = 100107
$ws.Range("H15").Value = 'Otros'
$ws.Range("I15").Value = 100107001
$ws.Range("J15").Value = 'Caqui'
$ws.Range("K15").Value = 'Mankaki'
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 130
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 14462
$ws.Range("Q15").Value = '$/caja 18 kilos granel'
$ws.Range("R15").Value = 'Región del Maule'
$ws.Range("S15").Value = 803
$ws.Range("T15").Value = 18

$ws.Range("A16").Value = 11
$ws.Range("B16").Value = 'Vega Monumental Concepción'
$ws.Range("C16").Value = 'Bíobío'
$ws.Range("D16").Value = 45091
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = 'Otros'
$ws.Range("I16").Value = 100107001
$ws.Range("J16").Value = 'Caqui'
$ws.Range("K16").Value = 'Mankaki'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 220
$ws.Range("N16").Value = 18000
$ws.Range("O16").Value = 19000
$ws.Range("P16").Value = 18455
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Provincia de Curicó'
$ws.Range("S16").Value = 1025
$ws.Range("T16").Value = 18

$ws.Range("A17").Value = 11
$ws.Range("B17").Value = 'Vega Monumental Concepción'
$ws.Range("C17").Value = 'Bíobío'
$ws.Range("D17").Value = 45091
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 'Fruta'
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = 'Otros'
$ws.Range("I17").Value = 100107001
$ws.Range("J17").Value = 'Caqui'
$ws.Range("K17").Value = 'Mankaki'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 15000
$ws.Range("O17").Value = 15000
$ws.Range("P17").Value = 15000
$ws.Range("Q17").Value = '$/caja 18 kilos granel'
$ws.Range("R17").Value = 'Provincia de Curicó'
$ws.Range("S17").Value = 833
$ws.Range("T17").Value = 18

$ws.Range("A18").Value = 11
$ws.Range("B18").Value = 'Vega Monumental Concepción'
$ws.Range("C18").Value = 'Bíobío'
$ws.Range("D18").Value = 45097
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 'Fruta'
$ws.Range("G18").Value = 100107
$ws.Range("H18").Value = 'Otros'
$ws.Range("I18").Value = 100107001
$ws.Range("J18").Value = 'Caqui'
$ws.Range("K18").Value = 'Mankaki'
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 18000
$ws.Range("O18").Value = 20000
$ws.Range("P18").Value = 19000
$ws.Range("Q18").Value = '$/caja 18 kilos granel'
$ws.Range("R18").Value = 'Región del Maule'
$ws.Range("S18").Value = 1056
$ws.Range("T18").Value = 18

$ws.Range("A19").Value = 11
$ws.Range("B19").Value = 'Vega Monumental Concepción'
$ws.Range("C19").Value = 'Bíobío'
$ws.Range("D19").Value = 45114
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = 'Fruta'
$ws.Range("G19").Value = 100107
$ws.Range("H19").Value = 'Otros'
$ws.Range("I19").Value = 100107001
$ws.Range("J19").Value = 'Caqui'
$ws.Range("K19").Value = 'Mankaki'
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Región del Maule'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18

$ws.Range("A20").Value = 11
$ws.Range("B20").Value = 'Vega Monumental Concepción'
$ws.Range("C20").Value = 'Bíobío'
$ws.Range("D20").Value = 45077
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 'Fruta'
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = 'Otros'
$ws.Range("I20").Value = 100107001
$ws.Range("J20").Value = 'Caqui'
$ws.Range("K20").Value = 'Mankaki'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = 12000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 12857
$ws.Range("Q20").Value = '$/caja 12 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 12857
$ws.Range("T20").Value = 1

$ws.Range("A21").Value = 11
$ws.Range("B21").Value = 'Vega Monumental Concepción'
$ws.Range("C21").Value = 'Bíobío'
$ws.Range("D21").Value = 45077
$ws.Range("E21").Value = 8
$ws.Range("F21").Value = 'Fruta'
$ws.Range("G21").Value = 100107
$ws.Range("H21").Value = 'Otros'
$ws.Range("I21").Value = 100107001
$ws.Range("J21").Value = 'Caqui'
$ws.Range("K21").Value = 'Mankaki'
$ws.Range("L21").Value = 'Segunda'
$ws.Range("M21").Value = 80
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 11000
$ws.Range("P21").Value = 11000
$ws.Range("Q21").Value = '$/caja 12 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 11000
$ws.Range("T21").Value = 1
